# This workbook lists weekly price observations (rows 2-19) that have been
# re-sorted/re-ordered (e.g. by date). The content of every row A:R is the
# same set of records, just relocated to different row numbers. We snapshot
# the current values first (so moves don't clobber each other), then write
# each snapshot back out at its new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the current (pre-edit) contents of every data row (2-19), columns A-R,
# before we start overwriting anything.
$firstRow = 2
$lastRow = 19
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $snapshot[$r] = $ws.Range("A" + $r + ":R" + $r).Value2
}

# Mapping: new row number -> old row number that supplies its data.
$mapping = @{
    2  = 5
    3  = 9
    4  = 2
    5  = 6
    6  = 8
    7  = 11
    8  = 14
    9  = 3
    10 = 7
    11 = 19
    12 = 18
    13 = 4
    14 = 10
    15 = 16
    16 = 17
    17 = 12
    18 = 13
    19 = 15
}

foreach ($newRow in $mapping.Keys) {
    $oldRow = $mapping[$newRow]
    $ws.Range("A" + $newRow + ":R" + $newRow).Value = $snapshot[$oldRow]
}
